{"js": "// Update the title date line (first paragraph of the document body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst datePara = paragraphs.items[0];\ndatePara.insertText(\"2025-06-26 Thursday\", Word.InsertLocation.replace);\n\n// Update every equation cell in the practice table. The table is laid out\n// as 20 rows x 5 columns, each cell holding a single run with the equation\n// text. Assigning the full 2-D `values` array rewrites each cell's text in\n// reading order while leaving cell/run formatting (fonts, size, alignment,\n// borders, etc.) untouched.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = [\n  [\"51+30=\", \"27-23=\", \"36+11=\", \"24-16=\", \"11+18=\"],\n  [\"15+19=\", \"0+51=\", \"75-51=\", \"24+9=\", \"27+60=\"],\n  [\"96-22=\", \"75+9=\", \"0+61=\", \"59+13=\", \"9+56=\"],\n  [\"12+82=\", \"41-8=\", \"43-26=\", \"11+13=\", \"5+67=\"],\n  [\"45-44=\", \"67-53=\", \"21+58=\", \"43-36=\", \"86-4=\"],\n  [\"86-25=\", \"80+7=\", \"74-33=\", \"2+78=\", \"44-44=\"],\n  [\"10+54=\", \"53+39=\", \"91-10=\", \"73+11=\", \"56-25=\"],\n  [\"49+48=\", \"15+7=\", \"27+68=\", \"36+16=\", \"17-7=\"],\n  [\"14+44=\", \"44-18=\", \"80-12=\", \"24+47=\", \"12+31=\"],\n  [\"84-59=\", \"7+74=\", \"73+8=\", \"50+33=\", \"72-61=\"],\n  [\"54+12=\", \"92-84=\", \"91-37=\", \"72+7=\", \"57-52=\"],\n  [\"64-60=\", \"35+5=\", \"87-42=\", \"89-27=\", \"74-57=\"],\n  [\"31+2=\", \"97-49=\", \"95-53=\", \"27-10=\", \"98-45=\"],\n  [\"60-33=\", \"14+72=\", \"59-37=\", \"79-64=\", \"66-40=\"],\n  [\"76-36=\", \"49+27=\", \"1+63=\", \"16+39=\", \"93-34=\"],\n  [\"61-59=\", \"57+4=\", \"35+54=\", \"32-12=\", \"77-4=\"],\n  [\"9+57=\", \"1+16=\", \"60-49=\", \"78-32=\", \"11+65=\"],\n  [\"15+2=\", \"2+52=\", \"4+29=\", \"30-6=\", \"9+35=\"],\n  [\"6+0=\", \"55-22=\", \"48-45=\", \"93-46=\", \"43-6=\"],\n  [\"89-60=\", \"3+15=\", \"42+28=\", \"41+58=\", \"38-30=\"]\n];\n\nawait context.sync();\n", "ps1": "# Update the title date line (first paragraph of the document).\n$d = $word.ActiveDocument\n$d.Paragraphs(1).Range.Text = \"2025-06-26 Thursday\"\n\n# Update every equation cell in the practice table. The table is laid out\n# as 20 rows x 5 columns, each cell holding a single run with the equation\n# text. Walk the grid in reading order and overwrite each cell's Range.Text\n# in place, which rewrites the text run while leaving cell/paragraph/run\n# formatting (fonts, size, alignment, borders, etc.) untouched.\n$newValues = @(\n    @(\"51+30=\", \"27-23=\", \"36+11=\", \"24-16=\", \"11+18=\"),\n    @(\"15+19=\", \"0+51=\", \"75-51=\", \"24+9=\", \"27+60=\"),\n    @(\"96-22=\", \"75+9=\", \"0+61=\", \"59+13=\", \"9+56=\"),\n    @(\"12+82=\", \"41-8=\", \"43-26=\", \"11+13=\", \"5+67=\"),\n    @(\"45-44=\", \"67-53=\", \"21+58=\", \"43-36=\", \"86-4=\"),\n    @(\"86-25=\", \"80+7=\", \"74-33=\", \"2+78=\", \"44-44=\"),\n    @(\"10+54=\", \"53+39=\", \"91-10=\", \"73+11=\", \"56-25=\"),\n    @(\"49+48=\", \"15+7=\", \"27+68=\", \"36+16=\", \"17-7=\"),\n    @(\"14+44=\", \"44-18=\", \"80-12=\", \"24+47=\", \"12+31=\"),\n    @(\"84-59=\", \"7+74=\", \"73+8=\", \"50+33=\", \"72-61=\"),\n    @(\"54+12=\", \"92-84=\", \"91-37=\", \"72+7=\", \"57-52=\"),\n    @(\"64-60=\", \"35+5=\", \"87-42=\", \"89-27=\", \"74-57=\"),\n    @(\"31+2=\", \"97-49=\", \"95-53=\", \"27-10=\", \"98-45=\"),\n    @(\"60-33=\", \"14+72=\", \"59-37=\", \"79-64=\", \"66-40=\"),\n    @(\"76-36=\", \"49+27=\", \"1+63=\", \"16+39=\", \"93-34=\"),\n    @(\"61-59=\", \"57+4=\", \"35+54=\", \"32-12=\", \"77-4=\"),\n    @(\"9+57=\", \"1+16=\", \"60-49=\", \"78-32=\", \"11+65=\"),\n    @(\"15+2=\", \"2+52=\", \"4+29=\", \"30-6=\", \"9+35=\"),\n    @(\"6+0=\", \"55-22=\", \"48-45=\", \"93-46=\", \"43-6=\"),\n    @(\"89-60=\", \"3+15=\", \"42+28=\", \"41+58=\", \"38-30=\")\n)\n\n$t = $d.Tables(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $t.Cell($r, $c).Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n"}
